$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "南京水軍左衞" / "Left Guard of Nanking Navy" office record that
# currently sits in row 2. The c_office_id column (A) keeps its original
# literal id values in rows 2-4 untouched; only the office name / dynasty /
# translation / pinyin / source columns (B:F) shift up from the rows below
# (Copy preserves each cell's original text/number typing, unlike a plain
# Value re-assignment which would coerce numeric-looking text like "19"
# into real numbers). The now-spare last row is then removed entirely.
$ws.Range("B3:F5").Copy($ws.Range("B2:F4"))
$ws.Rows(5).Delete()
